$wb = $excel.ActiveWorkbook

# Update the "Share of Cost Effective Capacity Built in a Single Year" sheet:
# the "onshore wind es" row (row 7) changes from 0.33 to 0.2 across all year columns (B:AE).
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$ws.Range("B7:AE7").Value = 0.2

# Move the active tab / selection to this sheet, selecting the row that was edited
# (mirrors the author re-saving with this sheet/range focused).
$ws.Activate()
$ws.Range("B7:AE7").Select()
